# Auto commit at 2025-09-27 16:15:13.42
# Append two new daily rows (2025-09-26, serial date 45926) for the two
# stations, mirroring the existing row layout, then move the selection
# to reflect the new "next empty" cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number formats (date / currency / integer styles) from the
# last existing pair of rows down onto the two new rows so the new
# cells reuse the same style indices instead of creating new ones.
$ws.Range("A50:F51").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 52: 四方坪站 station, 2025-09-26
$ws.Range("A52").Value = 45926
$ws.Range("B52").Value = "四方坪站"
$ws.Range("C52").Value = 8685.5400000000009
$ws.Range("D52").Value = 7050.51
$ws.Range("E52").Value = 3007.46
$ws.Range("F52").Value = 370

# Row 53: 高岭站 station, 2025-09-26
$ws.Range("A53").Value = 45926
$ws.Range("B53").Value = "高岭站"
$ws.Range("C53").Value = 3685.98
$ws.Range("D53").Value = 2938.61
$ws.Range("E53").Value = 1032.4100000000001
$ws.Range("F53").Value = 148

# Move the active selection to G52, matching the post-edit cursor
# position in the original workbook.
$ws.Range("G52").Select()
